$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

$ws.Range("A6").Value = "R-23-0140"
$ws.Range("B6").Value = "r230140@famt.ac.in"
$ws.Range("C6").Value = "Vaidehi Santosh Bhuwad"
$ws.Range("D6").Value = "saamia.kb@gmail.com"
$ws.Range("E6").Value = '$2b$12$3k654BgA1aTf6RyZ2ZZYTeV9x6kUM2lDTLfH/Clk7QCL0K0b9erae'
$ws.Range("F6").Value = "student"
